$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume (E) columns to text format so numeric-looking
# strings (e.g. "272.57") are not auto-converted into floating point numbers,
# matching the inlineStr text cells used in the workbook.
$priceRange = $ws.Range("D2:D51")
$volRange = $ws.Range("E2:E51")
$priceRange.NumberFormat = "@"
$volRange.NumberFormat = "@"

$ws.Range('D2').Value = '44.077.44'
$ws.Range('E2').Value = '  +1.66%  '
$ws.Range('D3').Value = '2.256.14'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '272.57'
$ws.Range('E5').Value = '  +5.48%  '
$ws.Range('D6').Value = '87.68'
$ws.Range('E6').Value = '  +10.78%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '0.613'
$ws.Range('E9').Value = '  +2.08%  '
$ws.Range('D10').Value = '45.63'
$ws.Range('E10').Value = '  +5.76%  '
$ws.Range('D11').Value = '0.0930'
$ws.Range('E11').Value = '  +0.74%  '
$ws.Range('D12').Value = '7.68'
$ws.Range('E12').Value = '  +8.39%  '
$ws.Range('D13').Value = '0.105'
$ws.Range('E13').Value = '  +2.25%  '
$ws.Range('D14').Value = '2.593.41'
$ws.Range('E14').Value = '  +1.06%  '
$ws.Range('D15').Value = '15.03'
$ws.Range('E15').Value = '  +2.63%  '
$ws.Range('D16').Value = '2.231.70'
$ws.Range('E16').Value = '  -0.11%  '
$ws.Range('E17').Value = '  +0.86%  '
$ws.Range('D18').Value = '44.013.70'
$ws.Range('E18').Value = '  +1.69%  '
$ws.Range('D19').Value = '0.0000103'
$ws.Range('E19').Value = '  -1.43%  '
$ws.Range('E20').Value = '  -0.27%  '
$ws.Range('D21').Value = '70.58'
$ws.Range('E21').Value = '  -1.05%  '
$ws.Range('D22').Value = '2.38'
$ws.Range('E22').Value = '  +2.52%  '
$ws.Range('D23').Value = '234.25'
$ws.Range('E23').Value = '  +1.83%  '
$ws.Range('D24').Value = '8.86'
$ws.Range('E24').Value = '  -4.98%  '
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('D26').Value = '2.52'
$ws.Range('E26').Value = '  +13.54%  '
$ws.Range('D27').Value = '10.89'
$ws.Range('E27').Value = '  +0.55%  '
$ws.Range('E28').Value = '  +6.22%  '
$ws.Range('D29').Value = '40.03'
$ws.Range('E29').Value = '  -4.33%  '
$ws.Range('E30').Value = '  +4.96%  '
$ws.Range('D31').Value = '175.54'
$ws.Range('E31').Value = '  +1.59%  '
$ws.Range('D32').Value = '20.94'
$ws.Range('E32').Value = '  +2.23%  '
$ws.Range('D33').Value = '0.0898'
$ws.Range('E33').Value = '  +3.35%  '
$ws.Range('E34').Value = '  +2.66%  '
$ws.Range('E35').Value = '  +1.40%  '
$ws.Range('D36').Value = '0.111'
$ws.Range('E36').Value = '  +3.49%  '
$ws.Range('D37').Value = '0.0352'
$ws.Range('E37').Value = '  -4.65%  '
$ws.Range('D38').Value = '4.38'
$ws.Range('E38').Value = '  -1.99%  '
$ws.Range('D39').Value = '3.51'
$ws.Range('E39').Value = '  +21.61%  '
$ws.Range('D40').Value = '12.72'
$ws.Range('E40').Value = '  -3.83%  '
$ws.Range('D41').Value = '2.22'
$ws.Range('E41').Value = '  +3.58%  '
$ws.Range('D42').Value = '64.66'
$ws.Range('E42').Value = '  +5.06%  '
$ws.Range('D43').Value = '5.46'
$ws.Range('E44').Value = '  +0.35%  '
$ws.Range('D45').Value = '8.51'
$ws.Range('E45').Value = '  -0.99%  '
$ws.Range('D46').Value = '0.0989'
$ws.Range('E46').Value = '  +0.75%  '
$ws.Range('D47').Value = '100.88'
$ws.Range('E47').Value = '  -2.58%  '
$ws.Range('D48').Value = '1.21'
$ws.Range('E48').Value = '  +4.74%  '
$ws.Range('E49').Value = '  +1.74%  '
$ws.Range('B50').Value = 'WOONetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D50').Value = '0.430'
$ws.Range('E50').Value = '  -8.76%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value = '1.49'
$ws.Range('E51').Value = '  +1.04%  '

# Restore the default (unstyled) cell style so we do not introduce any
# stray formatting differences versus the original workbook.
$priceRange.Style = "Normal"
$volRange.Style = "Normal"
